$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("en")
$ws2 = $wb.Worksheets.Item("de")

# --- Copy column-A formatting (style index 1, same as row 234) down into the
# new rows 235:239 on both sheets before writing any values, so the new A
# cells inherit the same style as the rest of the key column. On sheet "en",
# column B has a column-level default style (5) that new values pick up
# automatically; sheet "de" has no such column default, so its column-B
# cells need their format copied explicitly from an existing s="5" cell
# (B233) instead.
$ws1.Range("A234").Copy()
$ws1.Range("A235:A239").PasteSpecial(-4122)
$ws2.Range("A234").Copy()
$ws2.Range("A235:A239").PasteSpecial(-4122)
$ws2.Range("B233").Copy()
$ws2.Range("B235:B239").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Write new cell content in the exact order that reproduces the target
# shared-string table (new unique strings appended in this sequence):
#   Year, Years, Jahr, Jahre, Month, Monat, Months, Monate,
#   Additional info, Erweiterte Angaben, AboutInfo
$ws1.Range("A235").Value = "Year"
$ws1.Range("A236").Value = "Years"
$ws2.Range("B235").Value = "Jahr"
$ws2.Range("B236").Value = "Jahre"
$ws1.Range("A237").Value = "Month"
$ws2.Range("B237").Value = "Monat"
$ws1.Range("A238").Value = "Months"
$ws2.Range("B238").Value = "Monate"
$ws1.Range("B2").Value = "Additional info"
$ws2.Range("B2").Value = "Erweiterte Angaben"
$ws1.Range("A239").Value = "AboutInfo"

# --- Remaining cells simply reuse already-introduced shared strings.
$ws1.Range("B235").Value = "Year"
$ws1.Range("B236").Value = "Years"
$ws1.Range("B237").Value = "Month"
$ws1.Range("B238").Value = "Months"
$ws1.Range("B239").Value = "About"

$ws2.Range("A235").Value = "Year"
$ws2.Range("A236").Value = "Years"
$ws2.Range("A237").Value = "Month"
$ws2.Range("A238").Value = "Months"
$ws2.Range("A239").Value = "AboutInfo"
$ws2.Range("B239").Value = "Über"

# --- Selection / active-sheet state. "de" ends up the active tab, "en" keeps
# its own last selection from before it lost focus.
$ws1.Activate()
$ws1.Range("A239").Select()

$ws2.Activate()
$ws2.Range("B240").Select()
